$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.316871309442218
$ws.Range("D2").Value = 3.189181996700795
$ws.Range("E2").Value = 11.83991481257103
$ws.Range("F2").Value = 17.7689967141744
$ws.Range("G2").Value = 19.69113168747871
$ws.Range("H2").Value = 10.74234297618823
$ws.Range("I2").Value = 17.48216143757114
$ws.Range("M2").Value = 15.48507493782207
$ws.Range("N2").Value = 16.88662398178348
$ws.Range("O2").Value = 15.24879703180811
$ws.Range("B3").Value = 7.114463285498944
$ws.Range("D3").Value = 3.17571788092511
$ws.Range("E3").Value = 11.94500491597421
$ws.Range("F3").Value = 17.44253122203997
$ws.Range("G3").Value = 18.9383305960812
$ws.Range("H3").Value = 10.71860679258682
$ws.Range("I3").Value = 17.51951777260732
$ws.Range("M3").Value = 14.89619757031976
$ws.Range("N3").Value = 16.82698049310977
$ws.Range("O3").Value = 15.08326907121727
$ws.Range("B4").Value = 6.987625936183285
$ws.Range("D4").Value = 3.167868153985334
$ws.Range("E4").Value = 12.01527812415634
$ws.Range("F4").Value = 17.24495797617388
$ws.Range("G4").Value = 18.46764616574811
$ws.Range("H4").Value = 10.70651598236726
$ws.Range("I4").Value = 17.54625764350017
$ws.Range("M4").Value = 14.52400348487572
$ws.Range("N4").Value = 16.79229550405912
$ws.Range("O4").Value = 14.9854579227778
$ws.Range("B5").Value = 6.935367121982734
$ws.Range("D5").Value = 3.164776639672665
$ws.Range("E5").Value = 12.04534546919121
$ws.Range("F5").Value = 17.16529244964808
$ws.Range("G5").Value = 18.27406653570332
$ws.Range("H5").Value = 10.70221672634342
$ws.Range("I5").Value = 17.55811146265675
$ws.Range("M5").Value = 14.36988269699774
$ws.Range("N5").Value = 16.7786585692986
$ws.Range("O5").Value = 14.94660722110544
$ws.Range("B6").Value = 6.926657251637796
$ws.Range("D6").Value = 3.164269844842504
$ws.Range("E6").Value = 12.0504240961583
$ws.Range("F6").Value = 17.15211881655577
$ws.Range("G6").Value = 18.24182681649224
$ws.Range("H6").Value = 10.70154084606153
$ws.Range("I6").Value = 17.56013760762591
$ws.Range("M6").Value = 14.34415012707548
$ws.Range("N6").Value = 16.77642451387213
$ws.Range("O6").Value = 14.94021829829171
$ws.Range("B7").Value = 6.986923371608748
$ws.Range("D7").Value = 3.167826023165192
$ws.Range("E7").Value = 12.01567785050175
$ws.Range("F7").Value = 17.24387998602652
$ws.Range("G7").Value = 18.46504216292255
$ws.Range("H7").Value = 10.70645545495108
$ws.Range("I7").Value = 17.54641363200465
$ws.Range("M7").Value = 14.52193456150119
$ws.Range("N7").Value = 16.79210956366079
$ws.Range("O7").Value = 14.98492982684629
$ws.Range("B8").Value = 7.247653451505213
$ws.Range("D8").Value = 3.184454422200374
$ws.Range("E8").Value = 11.8749480201321
$ws.Range("F8").Value = 17.65590825575307
$ws.Range("G8").Value = 19.43352521717388
$ws.Range("H8").Value = 10.73364555009104
$ws.Range("I8").Value = 17.49425315342199
$ws.Range("M8").Value = 15.2843602431894
$ws.Range("N8").Value = 16.86566347569449
$ws.Range("O8").Value = 15.19095804453343
$ws.Range("B9").Value = 7.735745538417346
$ws.Range("D9").Value = 3.220264320616659
$ws.Range("E9").Value = 11.64525382453148
$ws.Range("F9").Value = 18.48131686033595
$ws.Range("G9").Value = 21.25164017988195
$ws.Range("H9").Value = 10.80649337930132
$ws.Range("I9").Value = 17.42209832231264
$ws.Range("M9").Value = 16.68663969492543
$ws.Range("N9").Value = 17.02478619747725
$ws.Range("O9").Value = 15.62316607782684
$ws.Range("B10").Value = 8.076670315651796
$ws.Range("D10").Value = 3.248372641749812
$ws.Range("E10").Value = 11.50564603489466
$ws.Range("F10").Value = 19.0912318315035
$ws.Range("G10").Value = 22.52163339266932
$ws.Range("H10").Value = 10.87164635224054
$ws.Range("I10").Value = 17.38738773103536
$ws.Range("M10").Value = 17.65039902128902
$ws.Range("N10").Value = 17.15010144464906
$ws.Range("O10").Value = 15.95489685016612
$ws.Range("B11").Value = 8.227280628616018
$ws.Range("D11").Value = 3.261515914274876
$ws.Range("E11").Value = 11.44867009302452
$ws.Range("F11").Value = 19.36800455965048
$ws.Range("G11").Value = 23.08227594161792
$ws.Range("H11").Value = 10.90374144186094
$ws.Range("I11").Value = 17.37555240760855
$ws.Range("M11").Value = 18.07271604467358
$ws.Range("N11").Value = 17.20878743611432
$ws.Range("O11").Value = 16.10823251808412
$ws.Range("B12").Value = 8.283621535047754
$ws.Range("D12").Value = 3.26654117368343
$ws.Range("E12").Value = 11.42805146839092
$ws.Range("F12").Value = 19.47259842506845
$ws.Range("G12").Value = 23.29192305043122
$ws.Range("H12").Value = 10.91624142863476
$ws.Range("I12").Value = 17.37163741810523
$ws.Range("M12").Value = 18.23020492707864
$ws.Range("N12").Value = 17.23123872448466
$ws.Range("O12").Value = 16.16659053563075
$ws.Range("B13").Value = 8.271519002218316
$ws.Range("D13").Value = 3.265456801069579
$ws.Range("E13").Value = 11.43244921662885
$ws.Range("F13").Value = 19.45008357150694
$ws.Range("G13").Value = 23.24689282431257
$ws.Range("H13").Value = 10.91353404929593
$ws.Range("I13").Value = 17.37245540555017
$ws.Range("M13").Value = 18.19639681326912
$ws.Range("N13").Value = 17.22639350447902
$ws.Range("O13").Value = 16.15400994873678
$ws.Range("B14").Value = 8.231929954973477
$ws.Range("D14").Value = 3.261928397633393
$ws.Range("E14").Value = 11.44695450675789
$ws.Range("F14").Value = 19.37661449375007
$ws.Range("G14").Value = 23.0995780170518
$ws.Range("H14").Value = 10.90476292543789
$ws.Range("I14").Value = 17.37521897138886
$ws.Range("M14").Value = 18.08572213762911
$ws.Range("N14").Value = 17.21063002134921
$ws.Range("O14").Value = 16.11302811405112
$ws.Range("B15").Value = 8.207589063223624
$ws.Range("D15").Value = 3.259773330934585
$ws.Range("E15").Value = 11.45596456097947
$ws.Range("F15").Value = 19.33158134470022
$ws.Range("G15").Value = 23.00899207163845
$ws.Range("H15").Value = 10.89943523702292
$ws.Range("I15").Value = 17.37698548807495
$ws.Range("M15").Value = 18.01761057528815
$ws.Range("N15").Value = 17.20100374827623
$ws.Range("O15").Value = 16.08796206468434
$ws.Range("B16").Value = 8.066733255632673
$ws.Range("D16").Value = 3.247520611442271
$ws.Range("E16").Value = 11.50950249320103
$ws.Range("F16").Value = 19.07312034253463
$ws.Range("G16").Value = 22.48463266668663
$ws.Range("H16").Value = 10.86959771173092
$ws.Range("I16").Value = 17.38824061960806
$ws.Range("M16").Value = 17.6224651415564
$ws.Range("N16").Value = 17.14629883881039
$ws.Range("O16").Value = 15.94492008265972
$ws.Range("B17").Value = 7.979140175604046
$ws.Range("D17").Value = 3.240093130123884
$ws.Range("E17").Value = 11.54403182323489
$ws.Range("F17").Value = 18.91430430752788
$ws.Range("G17").Value = 22.15842928067851
$ws.Range("H17").Value = 10.851917766588
$ws.Range("I17").Value = 17.39615689990527
$ws.Range("M17").Value = 17.37584295387903
$ws.Range("N17").Value = 17.11315990029458
$ws.Range("O17").Value = 15.85775070000464
$ws.Range("B18").Value = 7.92834001809814
$ws.Range("D18").Value = 3.235854799689438
$ws.Range("E18").Value = 11.56450597109979
$ws.Range("F18").Value = 18.82290008745323
$ws.Range("G18").Value = 21.96920788745274
$ws.Range("H18").Value = 10.84198041742533
$ws.Range("I18").Value = 17.401082497563
$ws.Range("M18").Value = 17.23248275947721
$ws.Range("N18").Value = 17.09425815289856
$ws.Range("O18").Value = 15.80784524003524
$ws.Range("B19").Value = 7.911069545430065
$ws.Range("D19").Value = 3.234425665495095
$ws.Range("E19").Value = 11.57154305211668
$ws.Range("F19").Value = 18.79194594734231
$ws.Range("G19").Value = 21.90487285506553
$ws.Range("H19").Value = 10.83865579856314
$ws.Range("I19").Value = 17.40281422773692
$ws.Range("M19").Value = 17.18368804303712
$ws.Range("N19").Value = 17.08788603787486
$ws.Range("O19").Value = 15.79098968204241
$ws.Range("B20").Value = 7.98850832934135
$ws.Range("D20").Value = 3.24088032636207
$ws.Range("E20").Value = 11.54029246525577
$ws.Range("F20").Value = 18.93121734783825
$ws.Range("G20").Value = 22.19332108543827
$ws.Range("H20").Value = 10.85377589705065
$ws.Range("I20").Value = 17.39527567318446
$ws.Range("M20").Value = 17.40225345435005
$ws.Range("N20").Value = 17.11667124769506
$ws.Range("O20").Value = 15.86700641989983
$ws.Range("B21").Value = 8.243577356548213
$ws.Range("D21").Value = 3.262963492493508
$ws.Range("E21").Value = 11.44266783966879
$ws.Range("F21").Value = 19.39820084803479
$ws.Range("G21").Value = 23.14292149251816
$ws.Range("H21").Value = 10.90732987777703
$ws.Range("I21").Value = 17.37439187879292
$ws.Range("M21").Value = 18.11829687527063
$ws.Range("N21").Value = 17.21525405326323
$ws.Range("O21").Value = 16.12505796416348
$ws.Range("B22").Value = 8.406227477560847
$ws.Range("D22").Value = 3.277675508410501
$ws.Range("E22").Value = 11.38445048662634
$ws.Range("F22").Value = 19.70210449780863
$ws.Range("G22").Value = 23.74798400932871
$ws.Range("H22").Value = 10.94434522721727
$ws.Range("I22").Value = 17.364045759998
$ws.Range("M22").Value = 18.57204514451246
$ws.Range("N22").Value = 17.28100725277728
$ws.Range("O22").Value = 16.29539227793583
$ws.Range("B23").Value = 8.319803488872859
$ws.Range("D23").Value = 3.269798905034684
$ws.Range("E23").Value = 11.41500521251405
$ws.Range("F23").Value = 19.54006101747893
$ws.Range("G23").Value = 23.42653345300872
$ws.Range("H23").Value = 10.92440753856746
$ws.Range("I23").Value = 17.36926617036488
$ws.Range("M23").Value = 18.33120786993675
$ws.Range("N23").Value = 17.24579688639435
$ws.Range("O23").Value = 16.20434619223693
$ws.Range("B24").Value = 7.984274359230964
$ws.Range("D24").Value = 3.240524335587286
$ws.Range("E24").Value = 11.54198109037941
$ws.Range("F24").Value = 18.92357126098311
$ws.Range("G24").Value = 22.17755172635383
$ws.Range("H24").Value = 10.85293512834662
$ws.Range("I24").Value = 17.39567290936553
$ws.Range("M24").Value = 17.39031816003133
$ws.Range("N24").Value = 17.11508329844723
$ws.Range("O24").Value = 15.86282125239364
$ws.Range("B25").Value = 7.606571762422146
$ws.Range("D25").Value = 3.210249440424133
$ws.Range("E25").Value = 11.70234217772378
$ws.Range("F25").Value = 18.25690375644544
$ws.Range("G25").Value = 20.77028412591103
$ws.Range("H25").Value = 10.78472132724196
$ws.Range("I25").Value = 17.43839844679339
$ws.Range("M25").Value = 16.31832871937795
$ws.Range("N25").Value = 16.98021316042236
$ws.Range("O25").Value = 15.50352933372446
